$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 holds the "Experimental" property; its value cell (B7) was empty and
# must now carry the literal text "true" (FHIR's required boolean rendered
# as text in this export, not an Excel TRUE/FALSE logical value).
# Writing the literal "true" directly would be auto-coerced to a Boolean by
# Excel, so build it as a text formula result first, then convert that
# formula to a plain text value in place (keeps the existing cell style).
$ws.Range("B7").Formula = '="true"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Row 8 holds the "Date" property; bump its value to the new publish date.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"

$excel.CutCopyMode = $false
